$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1:K1").Copy()
$ws.Range("A2:K2").PasteSpecial(-4122)

$ws.Range("A2").Value = "us-core-pregnancystatus"
$ws.Range("B2").Value = "US Core Pregnancy Status Observation Profile"
$ws.Range("C2").Value = "null#social-history"
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = "http://hl7.org/fhir/us/core/ValueSet/us-core-pregnancy-status-observation-codes (extensible)"
$ws.Range("G2").Value = "dateTime"
$ws.Range("H2").Value = "CodeableConcept"
$ws.Range("I2").Value = "optional"
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = ""
